# Auto-generated edit script: refresh market-price derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across multiple sheets.
$wb = $excel.ActiveWorkbook

# ==== Sheet: ALC ====
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3871.8125
$ws.Range("J40").Value = 4468.1816
$ws.Range("L40").Value = 4468.1816
$ws.Range("N40").Value = -4818.1816
# Row 42
$ws.Range("H42").Value = 292.33334
$ws.Range("I42").Value = 222.625
$ws.Range("K42").Value = 667.875
$ws.Range("M42").Value = -437.875
# Row 106
$ws.Range("H106").Value = 3456.5625
$ws.Range("I106").Value = 1258.7142
$ws.Range("K106").Value = 1258.7142
$ws.Range("M106").Value = -627.7141999999999
# Row 135
$ws.Range("H135").Value = 17475.666
$ws.Range("I135").Value = 927
$ws.Range("J135").Value = 25750
$ws.Range("K135").Value = 8343
$ws.Range("L135").Value = 231750
$ws.Range("M135").Value = -5808
$ws.Range("N135").Value = -236820
# Row 137
$ws.Range("H137").Value = 3346.9375
$ws.Range("I137").Value = 1309.2609
$ws.Range("K137").Value = 3927.7827
$ws.Range("M137").Value = -1377.7827
# Row 138
$ws.Range("H138").Value = 2546.4736
$ws.Range("I138").Value = 996.5
$ws.Range("J138").Value = 2728.8235
$ws.Range("K138").Value = 2989.5
$ws.Range("L138").Value = 8186.470499999999
$ws.Range("M138").Value = 2150.5
$ws.Range("N138").Value = -18466.4705

# ==== Sheet: ARM ====
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 5687811.5
$ws.Range("J74").Value = 19110.727
$ws.Range("L74").Value = 19110.727
$ws.Range("N74").Value = -20858.727
# Row 77
$ws.Range("H77").Value = 5687811.5
$ws.Range("J77").Value = 19110.727
$ws.Range("L77").Value = 95553.63499999999
$ws.Range("N77").Value = -104289.635
# Row 122
$ws.Range("H122").Value = 1227.5454
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
# Row 131
$ws.Range("H131").Value = 63199.8
$ws.Range("J131").Value = 63199.8
$ws.Range("L131").Value = 63199.8
$ws.Range("N131").Value = -73279.8
# Row 132
$ws.Range("H132").Value = 4793.1465
$ws.Range("I132").Value = 3122.1018
$ws.Range("K132").Value = 9366.305399999999
$ws.Range("M132").Value = -6836.305399999999

# ==== Sheet: BSM ====
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3090.3076
$ws.Range("I86").Value = 3319.4
$ws.Range("J86").Value = 2947.125
$ws.Range("K86").Value = 3319.4
$ws.Range("L86").Value = 2947.125
$ws.Range("M86").Value = -2196.4
$ws.Range("N86").Value = -5193.125
# Row 89
$ws.Range("H89").Value = 3090.3076
$ws.Range("I89").Value = 3319.4
$ws.Range("J89").Value = 2947.125
$ws.Range("K89").Value = 16597
$ws.Range("L89").Value = 14735.625
$ws.Range("M89").Value = -10981
$ws.Range("N89").Value = -25967.625
# Row 134
$ws.Range("H134").Value = 39720.414
$ws.Range("I134").Value = 1648.5
$ws.Range("K134").Value = 4945.5
$ws.Range("M134").Value = -2410.5

# ==== Sheet: CRP ====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 398088.38
$ws.Range("I31").Value = 2732.5405
$ws.Range("J31").Value = 1063005
$ws.Range("K31").Value = 2732.5405
$ws.Range("L31").Value = 1063005
$ws.Range("M31").Value = -2437.5405
$ws.Range("N31").Value = -1063595
# Row 34
$ws.Range("H34").Value = 398088.38
$ws.Range("I34").Value = 2732.5405
$ws.Range("J34").Value = 1063005
$ws.Range("K34").Value = 2732.5405
$ws.Range("L34").Value = 1063005
$ws.Range("M34").Value = -2530.5405
$ws.Range("N34").Value = -1063409
# Row 62
$ws.Range("H62").Value = 3997.5
$ws.Range("I62").Value = 3997.5
$ws.Range("K62").Value = 3997.5
$ws.Range("M62").Value = -3373.5
# Row 64
$ws.Range("H64").Value = 65000
$ws.Range("J64").Value = 65000
$ws.Range("L64").Value = 65000
$ws.Range("N64").Value = -65496
# Row 65
$ws.Range("H65").Value = 3997.5
$ws.Range("I65").Value = 3997.5
$ws.Range("K65").Value = 19987.5
$ws.Range("M65").Value = -16867.5
# Row 67
$ws.Range("H67").Value = 65000
$ws.Range("J67").Value = 65000
$ws.Range("L67").Value = 65000
$ws.Range("N67").Value = -66716
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# ==== Sheet: CUL ====
$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 14460.2
$ws.Range("I51").Value = 8384.571
$ws.Range("K51").Value = 25153.713
$ws.Range("M51").Value = -24693.713
# Row 112
$ws.Range("H112").Value = 9768.916999999999
$ws.Range("I112").Value = 9768.916999999999
$ws.Range("K112").Value = 29306.751
$ws.Range("M112").Value = -28198.751
# Row 134
$ws.Range("H134").Value = 3549.75
$ws.Range("I134").Value = 2152.625
$ws.Range("K134").Value = 6457.875
$ws.Range("M134").Value = -1387.875
# Row 136
$ws.Range("H136").Value = 6442.8335
$ws.Range("I136").Value = 6442.8335
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 19328.5005
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -14228.5005
$ws.Range("N136").ClearContents()
# Row 137
$ws.Range("H137").Value = 5071.091
$ws.Range("I137").Value = 7475
$ws.Range("J137").Value = 3697.4285
$ws.Range("K137").Value = 22425
$ws.Range("L137").Value = 11092.2855
$ws.Range("M137").Value = -17325
$ws.Range("N137").Value = -21292.2855
# Row 139
$ws.Range("H139").Value = 2324.0476
$ws.Range("I139").Value = 1972.1428
$ws.Range("K139").Value = 5916.428400000001
$ws.Range("M139").Value = -776.4284000000007
# Row 141
$ws.Range("H141").Value = 344776.56
$ws.Range("I141").Value = 1009329.7
$ws.Range("J141").Value = 12500
$ws.Range("K141").Value = 3027989.1
$ws.Range("L141").Value = 37500
$ws.Range("M141").Value = -3022809.1
$ws.Range("N141").Value = -47860

# ==== Sheet: GSM ====
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 7317.4287
$ws.Range("I102").Value = 4370.3335
$ws.Range("K102").Value = 4370.3335
$ws.Range("M102").Value = -2748.3335
# Row 126
$ws.Range("H126").Value = 4457
$ws.Range("J126").Value = 4457
$ws.Range("L126").Value = 13371
$ws.Range("N126").Value = -18311

# ==== Sheet: LTW ====
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 12567844
$ws.Range("I7").Value = 22226164
$ws.Range("J7").Value = 150002.86
$ws.Range("K7").Value = 22226164
$ws.Range("L7").Value = 150002.86
$ws.Range("M7").Value = -22226052
$ws.Range("N7").Value = -150226.86
# Row 40
$ws.Range("H40").Value = 3507.4614
$ws.Range("I40").Value = 1450
$ws.Range("J40").Value = 6799.4
$ws.Range("K40").Value = 1450
$ws.Range("L40").Value = 6799.4
$ws.Range("M40").Value = -1314
$ws.Range("N40").Value = -7071.4
# Row 46
$ws.Range("H46").Value = 4094
$ws.Range("I46").Value = 3673.1428
$ws.Range("J46").Value = 4935.7144
$ws.Range("K46").Value = 3673.1428
$ws.Range("L46").Value = 4935.7144
$ws.Range("M46").Value = -3485.1428
$ws.Range("N46").Value = -5311.7144
# Row 55
$ws.Range("H55").Value = 50000370
$ws.Range("I55").Value = 71428940
$ws.Range("K55").Value = 71428940
$ws.Range("M55").Value = -71428767
# Row 68
$ws.Range("H68").Value = 1332.8334
$ws.Range("I68").Value = 1399.6
$ws.Range("J68").Value = 999
$ws.Range("K68").Value = 1399.6
$ws.Range("L68").Value = 999
$ws.Range("M68").Value = -650.5999999999999
$ws.Range("N68").Value = -2497
# Row 71
$ws.Range("H71").Value = 1332.8334
$ws.Range("I71").Value = 1399.6
$ws.Range("J71").Value = 999
$ws.Range("K71").Value = 6998
$ws.Range("L71").Value = 4995
$ws.Range("M71").Value = -3254
$ws.Range("N71").Value = -12483
# Row 93
$ws.Range("H93").Value = 111113460
$ws.Range("I93").Value = 142859410
$ws.Range("J93").Value = 2604
$ws.Range("K93").Value = 142859410
$ws.Range("L93").Value = 2604
$ws.Range("M93").Value = -142858162
$ws.Range("N93").Value = -5100
# Row 126
$ws.Range("H126").Value = 12567844
$ws.Range("I126").Value = 22226164
$ws.Range("J126").Value = 150002.86
$ws.Range("K126").Value = 66678492
$ws.Range("L126").Value = 450008.58
$ws.Range("M126").Value = -66676022
$ws.Range("N126").Value = -454948.58
# Row 132
$ws.Range("H132").Value = 155634.22
$ws.Range("J132").Value = 194001
$ws.Range("L132").Value = 582003
$ws.Range("N132").Value = -587063

# ==== Sheet: WVR ====
$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 21666.334
$ws.Range("J40").Value = 21666.334
$ws.Range("L40").Value = 21666.334
$ws.Range("N40").Value = -21964.334
# Row 74
$ws.Range("H74").Value = 27500
$ws.Range("J74").Value = 27500
$ws.Range("L74").Value = 27500
$ws.Range("N74").Value = -29372
# Row 77
$ws.Range("H77").Value = 27500
$ws.Range("J77").Value = 27500
$ws.Range("L77").Value = 82500
$ws.Range("N77").Value = -91860
# Row 122
$ws.Range("H122").Value = 7810.3335
$ws.Range("I122").Value = 3505.875
$ws.Range("J122").Value = 11253.9
$ws.Range("K122").Value = 10517.625
$ws.Range("L122").Value = 33761.7
$ws.Range("M122").Value = -8067.625
$ws.Range("N122").Value = -38661.7
# Row 130
$ws.Range("H130").Value = 88745.5
$ws.Range("J130").Value = 88745.5
$ws.Range("L130").Value = 88745.5
$ws.Range("N130").Value = -98785.5
# Row 135
$ws.Range("H135").Value = 63331
$ws.Range("J135").Value = 63331
$ws.Range("L135").Value = 63331
$ws.Range("N135").Value = -73471
